$wb = $excel.ActiveWorkbook

# --- Astronauta (sheet 1) ---
$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Range("E2").Value = 1
$wsAstronauta.Range("F2").Value = 1
$wsAstronauta.Range("E3").Value = 1
$wsAstronauta.Range("F3").Value = 1
$wsAstronauta.Range("E4").Value = 1
$wsAstronauta.Range("F4").Value = 1
$wsAstronauta.Range("E5").Value = 1
$wsAstronauta.Range("E6").Value = 1
$wsAstronauta.Range("F6").Value = 1
$wsAstronauta.Range("E7").Value = 1
$wsAstronauta.Range("F7").Value = 1

# --- Mago (sheet 3) ---
$wsMago = $wb.Worksheets.Item("Mago")
$wsMago.Range("H3").Value = 1
$wsMago.Range("I4").Value = 1

# --- Ninja (sheet 4) ---
$wsNinja = $wb.Worksheets.Item("Ninja")
$wsNinja.Range("F2").Value = 1
$wsNinja.Range("F3").Value = 0
$wsNinja.Range("F4").Value = 1
$wsNinja.Range("F5").Value = 0
$wsNinja.Range("F6").Value = 1
$wsNinja.Range("F7").Value = 1

# --- Selections / active sheet ---
# Mago was previously the active (tabSelected) sheet with selection I7;
# it becomes non-active with selection I6.
$wsMago.Range("I6").Select()

# Ninja's selection moves from F4 to F3 (not the active tab).
$wsNinja.Range("F3").Select()

# Astronauta becomes the active tab, with selection F3.
$wsAstronauta.Activate()
$wsAstronauta.Range("F3").Select()
